$wb = $excel.ActiveWorkbook

# --- Sheet: DatosCuenta ---
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "PruebaCuenta"
$wsCuenta.Range("B2").Value = "Apellido"
$wsCuenta.Range("C2").Value = 27100101
$wsCuenta.Range("D2").Value = 106
$wsCuenta.Range("A4").Select()

# --- Sheet: DatosHogar ---
$wsHogar = $wb.Worksheets.Item("DatosHogar")
$wsHogar.Range("A2").Value = 627
$wsHogar.Range("A3").Select()

# --- Sheet: DatosMotor ---
$wsMotor = $wb.Worksheets.Item("DatosMotor")
$wsMotor.Range("A2").Value = "SMA008"
$wsMotor.Range("B2").Value = "ABC12SSMA008"
$wsMotor.Range("C2").Value = "ZAZ123SSMA008"
$wsMotor.Range("A2:C2").Select()

# --- Sheet: DatosAP ---
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Range("A2").Value = 21200107
$wsAP.Range("A2").Select()
